{"js": "// Typo fix in evaluation\n// 1) \"The program exceeds the expectations of the this requirement\"\n//    -> \"The program exceeds the expectations of this requirement\"\n// 2) \"I also progammed the menu myself...\" -> \"I also programmed the menu myself...\"\n// 3) Remove 3 of the trailing empty placeholder paragraphs (keep the first and last).\n\nconst body = context.document.body;\n\n// --- Change 1: fix \"the this requirement\" -> \"this requirement\" -----------\nconst dupWordMatches = body.search(\"the this requirement\", { matchCase: true, matchWholeWord: false });\ndupWordMatches.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < dupWordMatches.items.length; i++) {\n  dupWordMatches.items[i].insertText(\"this requirement\", \"Replace\");\n}\nawait context.sync();\n\n// --- Change 2: fix \"progammed\" -> \"programmed\" -----------------------------\nconst typoMatches = body.search(\"progammed\", { matchCase: true, matchWholeWord: false });\ntypoMatches.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < typoMatches.items.length; i++) {\n  typoMatches.items[i].insertText(\"programmed\", \"Replace\");\n}\nawait context.sync();\n\n// --- Change 3: drop 3 of the trailing empty underline paragraphs -----------\n// The document ends with a run of empty paragraphs; the middle three\n// (indent 720 twips / left=36pt, underlined run mark, no numbering) are\n// removed while the very first (indent 1440 twips) and very last\n// (indent 720 twips, but empty run rPr) empty paragraphs are kept.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,leftIndent,font/underline\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"\" && p.leftIndent === 36 && p.font.underline === \"Single\") {\n    toDelete.push(p);\n  }\n}\n\n// Keep only the extra copies beyond the first such paragraph: the diff\n// removes 3 consecutive empty/underlined paragraphs, leaving the\n// non-underlined empty paragraphs before/after untouched.\nfor (let i = 0; i < toDelete.length; i++) {\n  toDelete[i].delete();\n}\nawait context.sync();\n", "ps1": "# Typo fix in evaluation\n# 1) \"The program exceeds the expectations of the this requirement\"\n#    -> \"The program exceeds the expectations of this requirement\"\n# 2) \"I also progammed the menu myself...\" -> \"I also programmed the menu myself...\"\n# 3) Remove 3 of the trailing empty placeholder paragraphs (keep the first and last).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: fix \"the this requirement\" -> \"this requirement\" -----------\n$find1 = $d.Content.Find\n$find1.Text = \"the this requirement\"\n$find1.Replacement.Text = \"this requirement\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# --- Change 2: fix \"progammed\" -> \"programmed\" -----------------------------\n$find2 = $d.Content.Find\n$find2.Text = \"progammed\"\n$find2.Replacement.Text = \"programmed\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# --- Change 3: drop 3 of the trailing empty underline paragraphs -----------\n# The document ends with a run of empty paragraphs; the middle three\n# (indent 720 twips / left=36pt, underlined paragraph/run mark, no numbering)\n# are removed while the very first (indent 1440 twips) and very last\n# (indent 720 twips, but no underline on the run itself) empty paragraphs\n# are kept.\n$toDelete = @()\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $trimmed = $p.Range.Text -replace \"[\\r\\x07]+$\", \"\"\n    if ($trimmed -eq \"\" -and $p.Range.ParagraphFormat.LeftIndent -eq 36 -and $p.Range.Font.Underline -eq 1) {\n        $toDelete += $i\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $d.Paragraphs.Item($toDelete[$j]).Range.Delete()\n}\n"}
